$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.838.25'
$ws.Range("E2").Value = '  +0.69%  '
$ws.Range("D3").Value = '3.107.17'
$ws.Range("E3").Value = '  +0.71%  '
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").Value = '''575.45'
$ws.Range("E5").Value = '  -0.86%  '
$ws.Range("D6").Value = '''172.21'
$ws.Range("E6").Value = '  +2.91%  '
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("D8").Value = '3.105.27'
$ws.Range("E8").Value = '  +0.80%  '
$ws.Range("D9").Value = '''0.520'
$ws.Range("E9").Value = '  -0.45%  '
$ws.Range("D10").Value = '''6.41'
$ws.Range("E10").Value = '  -4.09%  '
$ws.Range("D11").Value = '''0.152'
$ws.Range("E11").Value = '  -0.78%  '
$ws.Range("D12").Value = '''0.477'
$ws.Range("E12").Value = '  -0.41%  '
$ws.Range("D13").Value = '''0.0000244'
$ws.Range("E13").Value = '  -1.91%  '
$ws.Range("D14").Value = '''36.93'
$ws.Range("E14").Value = '  +0.78%  '
$ws.Range("D15").Value = '''0.123'
$ws.Range("D16").Value = '3.628.83'
$ws.Range("E16").Value = '  +1.03%  '
$ws.Range("D17").Value = '66.935.01'
$ws.Range("E17").Value = '  +0.92%  '
$ws.Range("D18").Value = '''7.09'
$ws.Range("E18").Value = '  -1.18%  '
$ws.Range("D19").Value = '3.111.78'
$ws.Range("E19").Value = '  +0.92%  '
$ws.Range("D20").Value = '''16.35'
$ws.Range("E20").Value = '  +1.88%  '
$ws.Range("D21").Value = '''474.48'
$ws.Range("E21").Value = '  +2.45%  '
$ws.Range("D22").Value = '''0.709'
$ws.Range("E22").Value = '  -0.47%  '
$ws.Range("D23").Value = '''7.89'
$ws.Range("E23").Value = '  +6.08%  '
$ws.Range("D24").Value = '''13.39'
$ws.Range("E24").Value = '  +4.73%  '
$ws.Range("D25").Value = '''83.68'
$ws.Range("E25").Value = '  +0.65%  '
$ws.Range("D26").Value = '''2.27'
$ws.Range("E26").Value = '  -0.22%  '
$ws.Range("E27").Value = '  +0.02%  '
$ws.Range("D28").Value = '''9.86'
$ws.Range("E28").Value = '  -2.39%  '
$ws.Range("D29").Value = '''2.41'
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("D30").Value = '''7.88'
$ws.Range("E30").Value = '  -2.29%  '
$ws.Range("D31").Value = '''2.65'
$ws.Range("E31").Value = '  -0.43%  '
$ws.Range("D32").Value = '''28.54'
$ws.Range("E32").Value = '  +0.62%  '
$ws.Range("E33").Value = '  -1.28%  '
$ws.Range("D34").Value = '0.0₃0923'
$ws.Range("E34").Value = '  -9.49%  '
$ws.Range("E35").Value = '  -0.12%  '
$ws.Range("D36").Value = '''5.82'
$ws.Range("E36").Value = '  -1.03%  '
$ws.Range("D37").Value = '''0.977'
$ws.Range("E37").Value = '  -2.41%  '
$ws.Range("D38").Value = '''47.17'
$ws.Range("E38").Value = '  -2.67%  '
$ws.Range("D39").Value = '''2.06'
$ws.Range("E39").Value = '  +1.23%  '
$ws.Range("D40").Value = '''49.90'
$ws.Range("E40").Value = '  -0.42%  '
$ws.Range("D41").Value = '''0.308'
$ws.Range("E41").Value = '  -1.77%  '
$ws.Range("E42").Value = '  -0.32%  '
$ws.Range("D43").Value = '''8.57'
$ws.Range("E43").Value = '  -1.01%  '
$ws.Range("D44").Value = '2.788.78'
$ws.Range("E44").Value = '  +0.49%  '
$ws.Range("D45").Value = '''0.0354'
$ws.Range("E45").Value = '  -1.83%  '
$ws.Range("D46").Value = '''377.11'
$ws.Range("E46").Value = '  -2.27%  '
$ws.Range("D47").Value = '''2.52'
$ws.Range("E47").Value = '  -12.95%  '
$ws.Range("D48").Value = '''135.66'
$ws.Range("E48").Value = '  +0.85%  '
$ws.Range("E49").Value = '  -0.04%  '
$ws.Range("D50").Value = '''24.77'
$ws.Range("E50").Value = '  +1.08%  '
$ws.Range("D51").Value = '''2.19'
$ws.Range("E51").Value = '  -1.31%  '